# Auto-generated Excel COM-interop script
# Applies the BRVM automatic update (GitHub Actions) to both sheets:
#  - 'Recommandations': rows 2..31 fully refreshed (row 32 removed)
#  - 'Top_YTD': rows 2..11 values refreshed

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd  = $wb.Worksheets.Item("Top_YTD")

# ---- Target data for 'Recommandations' (A:G, rows 2..31) ----
$recoData = @(
    @('BRVM - CONSOMMATION DE BASE     (**)', 0, 4, 1108.45, 284.89, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRINCIPAL     (**)', 0, 4, 1054.73, 269.95, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 4, 757.59, 190.02, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIELS', 0, 4, 755.99, 193.58, '🟡 Observer', '➖ Neutre'),
    @('BRVM - SERVICES FINANCIERS', 0, 4, 650.95, 164.6, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRESTIGE', 0, 4, 618.87, 154.99, '🟡 Observer', '➖ Neutre'),
    @('BRVM – COMPOSITE TOTAL RETURN     (**)', 0, 4, 598.85, 151.44, '🟡 Observer', '➖ Neutre'),
    @('BRVM - ENERGIE', 0, 4, 506.29, 126.39, '🟡 Observer', '➖ Neutre'),
    @('BRVM - SERVICES PUBLICS', 0, 4, 478.48, 121.62, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TELECOMMUNICATIONS', 0, 4, 403.64, 101.42, '🟡 Observer', '➖ Neutre'),
    @('ECOBANK TRANS. INCORP. TG (ETIT)', 4, 0, 24.61, 6.45, '🟢 Achat', '✅ Renforcer'),
    @('UNILEVER CI (UNLC)', 3, 0, 19.83, 6.88, '🟢 Achat', '✅ Renforcer'),
    @('SUCRIVOIRE (SCRC)', 2, 0, 14.83, 7.49, '🟡 Observer', '➖ Neutre'),
    @('SICOR CI (SICC)', 1, 0, 7.44, 7.44, '🟡 Observer', '➖ Neutre'),
    @('NEI-CEDA CI (NEIC)', 1, 0, 6.92, 6.92, '🟡 Observer', '➖ Neutre'),
    @('CORIS BANK INTERNATIONAL (CBIBF)', 1, 1, 5.95, -1.55, '🟡 Observer', '👀 À surveiller'),
    @('SITAB CI (STBC)', 1, 1, 3.39, -1.9, '🟡 Observer', '👀 À surveiller'),
    @('SOLIBRA CI (SLBC)', 1, 0, 2.17, 2.17, '🟡 Observer', '➖ Neutre'),
    @('LOTERIE NATIONALE DU BENIN (LNBB)', 1, 0, 1.65, 1.65, '🟡 Observer', '➖ Neutre'),
    @('SMB CI (SMBC)', 1, 1, 1.13, -6.36, '🟡 Observer', '👀 À surveiller'),
    @('BANK OF AFRICA CI (BOAC)', 0, 1, -1.32, -1.32, '🟡 Observer', '➖ Neutre'),
    @('ONATEL BF (ONTBF)', 0, 1, -2.95, -2.95, '🟡 Observer', '➖ Neutre'),
    @('UNIWAX CI (UNXC)', 1, 1, -3.19, 4.27, '🟡 Observer', '👀 À surveiller'),
    @('ECOBANK COTE D''''IVOIRE (ECOC)', 1, 1, -3.54, -5.33, '🟡 Observer', '👀 À surveiller'),
    @('TRACTAFRIC MOTORS CI (PRSC)', 0, 1, -4.67, -4.67, '🟡 Observer', '➖ Neutre'),
    @('BERNABE CI (BNBC)', 0, 2, -4.79, -0.66, '🟡 Observer', '➖ Neutre'),
    @('ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)', 1, 2, -7.31, 7.45, '🟡 Observer', '👀 À surveiller'),
    @('SETAO CI (STAC)', 0, 2, -9.82, -2.99, '🟡 Observer', '➖ Neutre'),
    @('EVIOSYS PACKAGING SIEM CI (SEMC)', 1, 3, -14.62, 7.41, '🔴 Vente', '⚠️ Risque de décrochage'),
    @('SICABLE CI (CABC)', 0, 3, -22.38, -7.45, '🔴 Vente', '⚠️ Risque de décrochage')
)

# ---- Target data for 'Top_YTD' (A:B, rows 2..11) ----
$ytdData = @(
    @('BRVM - CONSOMMATION DE BASE     (**)', 20111.26),
    @('BRVM-PRINCIPAL     (**)', 17387.49),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 6914),
    @('BRVM - INDUSTRIELS', 6873.06),
    @('BRVM - SERVICES FINANCIERS', 4664.98),
    @('BRVM-PRESTIGE', 4109.51),
    @('BRVM – COMPOSITE TOTAL RETURN     (**)', 3788.06),
    @('BRVM - ENERGIE', 2535.21),
    @('BRVM - SERVICES PUBLICS', 2226.24),
    @('BRVM - TELECOMMUNICATIONS', 1529.29)
)

# Write the refreshed 'Recommandations' rows 2..31
$r = 2
foreach ($row in $recoData) {
    $wsReco.Cells.Item($r, 1).Value = $row[0]
    $wsReco.Cells.Item($r, 2).Value = $row[1]
    $wsReco.Cells.Item($r, 3).Value = $row[2]
    $wsReco.Cells.Item($r, 4).Value = $row[3]
    $wsReco.Cells.Item($r, 5).Value = $row[4]
    $wsReco.Cells.Item($r, 6).Value = $row[5]
    $wsReco.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# Old table had 31 data rows (r=2..32); new table only has 30 (r=2..31),
# so the trailing row 32 must be removed entirely (shift cells up).
$wsReco.Rows.Item(32).Delete()

# Write the refreshed 'Top_YTD' rows 2..11
$r = 2
foreach ($row in $ytdData) {
    $wsYtd.Cells.Item($r, 1).Value = $row[0]
    $wsYtd.Cells.Item($r, 2).Value = $row[1]
    $r++
}

